# Update "想去人数" (want-to-go count) values in column F
# for the "展览" (Exhibitions) and "全部类型" (All Types) sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 888
    $ws.Range("F3").Value = 4502
    $ws.Range("F4").Value = 128
    $ws.Range("F5").Value = 790
}
